# Update meeting minutes workbook per commit "[PHONGDV] Update meetting minus"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Week 1_Date 12_5"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 20.7109375
$ws.Columns.Item(2).ColumnWidth = 27.7109375
$ws.Columns.Item(3).ColumnWidth = 27.85546875
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 39

# Header row
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Location"
$ws.Range("C1").Value = "Item"
$ws.Range("D1").Value = "Who"
$ws.Range("E1").Value = "Action Item"
$ws.Range("F1").Value = "Note"

$headerRange = $ws.Range("A1:F1")
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.TintAndShade = -0.249977111117893

# Row 2
$ws.Range("A2").Value = 43804.375
$ws.Range("A2").NumberFormat = "m/d/yyyy h:mm"
$ws.Range("B2").Value = "Coffee house ấp bắc"
$ws.Range("C2").Value = "Identify Basic Requirement"
$ws.Range("D2").Value = "all team member"
$ws.Range("E2").Value = "Create excel sheet for basic story"

# Row 3
$ws.Range("C3").Value = "Decide technical stack"
$ws.Range("D3").Value = "all team member"
$ws.Range("E3").Value = "create technical stack sheet"

# Row 4
$ws.Range("C4").Value = "Raising github restructure"
$ws.Range("D4").Value = "phong"
$ws.Range("E4").Value = "restructure  git branch and permission"

# Selection matches recorded state: B4 active cell
$ws.Range("B4").Select()
